$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.181.64"
$ws.Range("E2").Value = "  +3.50%  "
$ws.Range("D3").Value = "3.209.03"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.15"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.26"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.19%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.537"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.34%  "
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("E10").Value = "  +4.19%  "
$ws.Range("E11").Value = "  +3.72%  "
$ws.Range("D12").Value = "3.764.74"
$ws.Range("E12").Value = "  +2.31%  "
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.18"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000174"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.27%  "
$ws.Range("D16").Value = "60.239.21"
$ws.Range("E16").Value = "  +3.53%  "
$ws.Range("D17").Value = "3.211.88"
$ws.Range("E17").Value = "  +2.42%  "
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.16"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.35"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.13"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.529"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.26"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.91"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +10.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.170"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.87%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "0.0₃0909"
$ws.Range("E28").Value = "  +3.12%  "
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.42"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.45"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.18"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  +3.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.62"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "157.09"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.73%  "
$ws.Range("E36").Value = "  +1.10%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "2.794.99"
$ws.Range("E37").Value = "  +5.62%  "
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.88"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0708"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.68"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.25"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.12"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.722"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0287"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.89%  "
$ws.Range("D45").Value = "3.253.08"
$ws.Range("E45").Value = "  +2.20%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.104"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.18"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.807"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +7.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.82"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.69%  "
$ws.Range("E51").Value = "  -0.01%  "
